$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (legacy password hash) which blocks direct
# cell writes via COM; unprotect first so the value updates below succeed.
$ws.Unprotect()

# Update the confidential notice date from 2021-05-14 to 2021-05-17
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05672176165359961
$ws.Range("E2").Value = -0.002134174631158836

$ws.Range("D3").Value = 0.02372200705140393
$ws.Range("E3").Value = -0.001553096486119232

$ws.Range("D4").Value = 0.03067143602452925
$ws.Range("E4").Value = 0.0005791505791505891

$ws.Range("D5").Value = 0.03326291578644733
$ws.Range("E5").Value = 0.01766253288237496

$ws.Range("D6").Value = 0.0381004762559532
$ws.Range("E6").Value = 0.007045171985080811

$ws.Range("D7").Value = 0.01944366409843281
$ws.Range("E7").Value = -0.004500236854571016

$ws.Range("D8").Value = 0.004263211184876654
$ws.Range("E8").Value = -0.009259259259259633

$ws.Range("D9").Value = 0.006947455264243434
$ws.Range("E9").Value = -0.006439393939393856

$ws.Range("D10").Value = 0.07231669343235211
$ws.Range("E10").Value = 0.01200873362445409

$ws.Range("D11").Value = 0.07243511596526536
$ws.Range("E11").Value = 0.01198910081743865

$ws.Range("D12").Value = 0.1442965405488622
$ws.Range("E12").Value = -0.002115552961774347

$ws.Range("D13").Value = 0.3830080770746477
$ws.Range("E13").Value = -0.0009636443276391793

$ws.Range("D14").Value = 0.1148106456593865
$ws.Range("E14").Value = 0.00144404332129966

$ws.Range("E15").Value = 0.001772390575934812
